$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.763.54"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.271.75"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'249.26"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'0.642"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("D7").Value = "'76.90"
$ws.Range("E7").Value = "  +7.19%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.639"
$ws.Range("E9").Value = "  -3.58%  "
$ws.Range("D10").Value = "'39.87"
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("D11").Value = "'0.0968"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "'7.28"
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "2.609.88"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "'0.865"
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "2.273.32"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "42.625.38"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "0.0₃0988"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "'6.17"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").Value = "'72.01"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").Value = "'235.16"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  -6.27%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'11.27"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").Value = "'167.39"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'20.86"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "'6.38"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").Value = "'0.0851"
$ws.Range("E32").Value = "  +5.63%  "
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("D34").Value = "'30.31"
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  -1.72%  "
$ws.Range("D38").Value = "'0.0305"
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").Value = "'13.69"
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("D40").Value = "'2.26"
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").Value = "'5.84"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'0.208"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "'108.98"
$ws.Range("E43").Value = "  +14.29%  "
$ws.Range("D44").Value = "'60.78"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("D46").Value = "'0.100"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("E47").Value = "  -8.59%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("D51").Value = "'4.19"
$ws.Range("E51").Value = "  -1.97%  "
